$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contextual")

# Insert a new column B ("country") between "community" (A) and "e_coli" (old B).
$ws.Columns.Item(2).Insert()

# Header
$ws.Range("B1").Value = "country"

# New country values per community/karst row, matching the region -> country mapping
$ws.Range("B2").Value = "Indonesia"
$ws.Range("B3").Value = "Philippines"
$ws.Range("B4").Value = "Bolivia"
$ws.Range("B5").Value = "Colombia"

# Nudge number format so these new cells pick up their own style entry
$ws.Range("B2:B5").NumberFormat = "General"

# Update selection to reflect where the user ended up after editing
$ws.Activate()
$ws.Range("H6").Select()
